# html QC report rendering; update README
#
# 1) The ".pdf" QC-report label becomes ".html" and its textbox is
#    nudged right (same y, new x) to match the new (shorter-looking)
#    rendering.
# 2) The cached "datetimeFigureOut" field text on every slide layout and
#    the slide master is refreshed from 10/14/16 -> 11/6/16 (this is the
#    auto-updating "today" date stamp that PowerPoint bakes into the
#    footer placeholder whenever the deck is saved on a different day).

$p = $ppt.ActivePresentation

# --- 1) Quality-control report textbox: .pdf -> .html -------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 44") {
        $shape.Left = 530.78
        $shape.TextFrame.TextRange.Text = ".html"
    }
}

# --- 2) Refresh the cached date stamp everywhere it is cached -----------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = "11/6/16"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $layout.Shapes
}
